$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in header label (shared string): "Clava_representante" -> "Clave_representante"
$ws.Range("B1").Value = "Clave_representante"

# Unify the formatting of K1:L1 with the rest of the header row (A1), which collapses the
# two near-identical header styles into a single one (K1/L1 previously carried a
# redundant style that only differed by an inert applyAlignment/applyProtection flag).
$ws.Range("A1").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122)

# Leave the header row selected (A1:M1, with A1 as the active cell) as the final
# selection state of the sheet.
$ws.Range("A1:M1").Select()
